$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 7374.625
$ws.Range("I4").Value = 7374.625
$ws.Range("K4").Value = 7374.625
$ws.Range("M4").Value = -7260.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1189.2222
$ws.Range("I28").Value = 1014.8571
$ws.Range("K28").Value = 1014.8571
$ws.Range("M28").Value = -529.8570999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2792.7576
$ws.Range("I33").Value = 3103.875
$ws.Range("J33").Value = 1963.1111
$ws.Range("K33").Value = 3103.875
$ws.Range("L33").Value = 1963.1111
$ws.Range("M33").Value = -2874.875
$ws.Range("N33").Value = -2421.1111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1478.4667
$ws.Range("J41").Value = 1248.625
$ws.Range("L41").Value = 1248.625
$ws.Range("N41").Value = -2128.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4916.5
$ws.Range("J64").Value = 4959.8
$ws.Range("L64").Value = 4959.8
$ws.Range("N64").Value = -5455.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4916.5
$ws.Range("J67").Value = 4959.8
$ws.Range("L67").Value = 4959.8
$ws.Range("N67").Value = -6675.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 70011.664
$ws.Range("J68").Value = 70011.664
$ws.Range("L68").Value = 70011.664
$ws.Range("N68").Value = -71509.664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H71").Value = 70011.664
$ws.Range("J71").Value = 70011.664
$ws.Range("L71").Value = 210034.992
$ws.Range("N71").Value = -217522.992

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3927.4443
$ws.Range("I76").Value = 3492.4285
$ws.Range("J76").Value = 5450
$ws.Range("K76").Value = 3492.4285
$ws.Range("L76").Value = 5450
$ws.Range("M76").Value = -3177.4285
$ws.Range("N76").Value = -6080

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3927.4443
$ws.Range("I79").Value = 3492.4285
$ws.Range("J79").Value = 5450
$ws.Range("K79").Value = 3492.4285
$ws.Range("L79").Value = 5450
$ws.Range("M79").Value = -2400.4285
$ws.Range("N79").Value = -7634

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 736018.4
$ws.Range("I88").Value = 1900
$ws.Range("K88").Value = 1900
$ws.Range("M88").Value = -1494

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 736018.4
$ws.Range("I91").Value = 1900
$ws.Range("K91").Value = 1900
$ws.Range("M91").Value = -496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3448.2
$ws.Range("I111").Value = 3387
$ws.Range("K111").Value = 10161
$ws.Range("M111").Value = -7094

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2215.238
$ws.Range("J132").Value = 2332.6667
$ws.Range("L132").Value = 6998.000100000001
$ws.Range("N132").Value = -12058.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3132.8975
$ws.Range("I137").Value = 1473.4231
$ws.Range("J137").Value = 6451.846
$ws.Range("K137").Value = 4420.2693
$ws.Range("L137").Value = 19355.538
$ws.Range("M137").Value = -1870.2693
$ws.Range("N137").Value = -24455.538

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2498.5386
$ws.Range("I141").Value = 2498.5386
$ws.Range("K141").Value = 7495.6158
$ws.Range("M141").Value = -2315.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4369.9
$ws.Range("I63").Value = 3924.125
$ws.Range("K63").Value = 3924.125
$ws.Range("M63").Value = -3238.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4369.9
$ws.Range("I66").Value = 3924.125
$ws.Range("K66").Value = 19620.625
$ws.Range("M66").Value = -16188.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1570.875
$ws.Range("I110").Value = 1570.875
$ws.Range("K110").Value = 1570.875
$ws.Range("M110").Value = 474.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 999
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2997
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -547
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5749.5127
$ws.Range("I132").Value = 2594.8965
$ws.Range("J132").Value = 14897.9
$ws.Range("K132").Value = 7784.689499999999
$ws.Range("L132").Value = 44693.7
$ws.Range("M132").Value = -5254.689499999999
$ws.Range("N132").Value = -49753.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3003.524
$ws.Range("J20").Value = 3341.6
$ws.Range("L20").Value = 3341.6
$ws.Range("N20").Value = -3835.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 816.4545000000001
$ws.Range("I22").Value = 465.125
$ws.Range("J22").Value = 1753.3334
$ws.Range("K22").Value = 465.125
$ws.Range("L22").Value = 1753.3334
$ws.Range("M22").Value = -292.125
$ws.Range("N22").Value = -2099.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 119999
$ws.Range("J63").Value = 119999
$ws.Range("L63").Value = 119999
$ws.Range("N63").Value = -121371

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H66").Value = 119999
$ws.Range("J66").Value = 119999
$ws.Range("L66").Value = 359997
$ws.Range("N66").Value = -366861

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 49840.5
$ws.Range("I134").Value = 890.9231
$ws.Range("K134").Value = 2672.7693
$ws.Range("M134").Value = -137.7692999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 525024.8
$ws.Range("I31").Value = 8996.071
$ws.Range("K31").Value = 8996.071
$ws.Range("M31").Value = -8701.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 525024.8
$ws.Range("I34").Value = 8996.071
$ws.Range("K34").Value = 8996.071
$ws.Range("M34").Value = -8794.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1257.4546
$ws.Range("J58").Value = 1253.3334
$ws.Range("L58").Value = 1253.3334
$ws.Range("N58").Value = -1659.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2849.125
$ws.Range("I62").Value = 2598
$ws.Range("J62").Value = 2999.8
$ws.Range("K62").Value = 2598
$ws.Range("L62").Value = 2999.8
$ws.Range("M62").Value = -1974
$ws.Range("N62").Value = -4247.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2849.125
$ws.Range("I65").Value = 2598
$ws.Range("J65").Value = 2999.8
$ws.Range("K65").Value = 12990
$ws.Range("L65").Value = 14999
$ws.Range("M65").Value = -9870
$ws.Range("N65").Value = -21239

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 437.25
$ws.Range("I107").Value = 307.375
$ws.Range("K107").Value = 307.375
$ws.Range("M107").Value = 1612.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1568.1072
$ws.Range("I132").Value = 1441
$ws.Range("K132").Value = 4323
$ws.Range("M132").Value = -1793

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1257.4546
$ws.Range("J136").Value = 1253.3334
$ws.Range("L136").Value = 3760.0002
$ws.Range("N136").Value = -8860.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1232.1666
$ws.Range("J117").Value = 594
$ws.Range("L117").Value = 1782
$ws.Range("N117").Value = -8666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 29668.666
$ws.Range("I5").Value = 20001
$ws.Range("J5").Value = 34502.5
$ws.Range("K5").Value = 20001
$ws.Range("L5").Value = 34502.5
$ws.Range("M5").Value = -19889
$ws.Range("N5").Value = -34726.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 142876290
$ws.Range("I132").Value = 166672340
$ws.Range("K132").Value = 500017020
$ws.Range("M132").Value = -500014490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5853.304
$ws.Range("I46").Value = 3545.1333
$ws.Range("K46").Value = 3545.1333
$ws.Range("M46").Value = -3357.1333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3165.6667
$ws.Range("I68").Value = 2998.5
$ws.Range("K68").Value = 2998.5
$ws.Range("M68").Value = -2249.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3165.6667
$ws.Range("I71").Value = 2998.5
$ws.Range("K71").Value = 14992.5
$ws.Range("M71").Value = -11248.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 29990
$ws.Range("I75").Value = 29990
$ws.Range("K75").Value = 29990
$ws.Range("M75").Value = -29054

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 29990
$ws.Range("I78").Value = 29990
$ws.Range("K78").Value = 89970
$ws.Range("M78").Value = -85290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 32268.4
$ws.Range("I136").Value = 4700.0386
$ws.Range("K136").Value = 14100.1158
$ws.Range("M136").Value = -11550.1158

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25004624
$ws.Range("I62").Value = 4999.5
$ws.Range("J62").Value = 33337834
$ws.Range("K62").Value = 4999.5
$ws.Range("L62").Value = 33337834
$ws.Range("M62").Value = -4375.5
$ws.Range("N62").Value = -33339082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 25004624
$ws.Range("I65").Value = 4999.5
$ws.Range("J65").Value = 33337834
$ws.Range("K65").Value = 24997.5
$ws.Range("L65").Value = 166689170
$ws.Range("M65").Value = -21877.5
$ws.Range("N65").Value = -166695410

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 31254072
$ws.Range("I107").Value = 41668896
$ws.Range("K107").Value = 125006688
$ws.Range("M107").Value = -125004768

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 886.26086
$ws.Range("I113").Value = 919.0476
$ws.Range("K113").Value = 2757.1428
$ws.Range("M113").Value = -587.1428000000001
